$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5: SMV Powertrain, Battery -> SMV Powertrain, Battery, CFD
$ws.Range("A5").Value = "SMV Powertrain, Battery, CFD"
$ws.Range("B5").Value = "images, CAD, etc, brochure, ansys results, images"

# Row 6 used to be "SMV CFD" / "ansys results, images" -- now becomes "TEG" / poster board text
$ws.Range("A6").Value = "TEG"
$ws.Range("B6").Value = "poster board, images, block diagram, analysis results"

# Row 7 used to be "TEG" / poster board text -- now becomes "BARC 131" / Drifting text
$ws.Range("A7").Value = "BARC 131"
$ws.Range("B7").Value = "Drifting, images, videos, simulink block diagram, code push to github"

# Row 8 used to be "BARC, QUAD" / images videos text -- now becomes "Quadcopter 136" / Code images videos
$ws.Range("A8").Value = "Quadcopter 136"
$ws.Range("B8").Value = "Code, images, videos"

# Row 12: Movi Pro highlights updated
$ws.Range("B12").Value = "final video, CAD, link to movipro"

# Row 16: add new project name, keep existing "Maybe some cool analysis project" text
$ws.Range("A16").Value = "C180 FEM Projects"

# Row 17: new row with only column A
$ws.Range("A17").Value = "E7 Marching Project"

# Update selection to match the diff (activeCell B12)
$ws.Range("B12").Select()
